# Juno: check in to OLPRODLOC.
#
# This script bolds several heading/job-title runs in the resume and
# renames the "动画团队经理" job title to "动画团队管理人员".

$d = $word.ActiveDocument

# 1) Name heading "Nestor Wilke" -> bold
$d.Paragraphs.Item(1).Range.Font.Bold = 1

# 2) "工作经历" (Work Experience) heading -> bold
$d.Paragraphs.Item(5).Range.Font.Bold = 1

# 3) "动画团队经理" job title -> bold, and text changes to "动画团队管理人员"
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Font.Bold = 1
$p6.Range.Find.Execute("动画团队经理", $false, $false, $false, $false, $false, `
                        $true, 1, $false, "动画团队管理人员", 2)

# 4) "高级动画设计师" (Senior Animation Designer) job title -> bold
$d.Paragraphs.Item(12).Range.Font.Bold = 1

# 5) "动画设计师" (Animation Designer) job title -> bold
$d.Paragraphs.Item(18).Range.Font.Bold = 1

# 6) "动画美术学士学位" (Bachelor's degree) -> bold
$d.Paragraphs.Item(24).Range.Font.Bold = 1
